# Insert a new price-observation row before the current row 30 ("Fruta,
# Terminal Hortofrutícola Agro Chillán - Mango" weekly/daily price series).
# All existing rows from 30 down to 149 shift down by one (to 31..150),
# and the freshly inserted row 30 is populated with a new observation.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a blank row at position 30; rows 30:149 shift to 31:150.
$ws.Rows(30).Insert()

# Copy the formatting/values of the row now sitting at 31 (the old row 30)
# into the new blank row 30, limited to the used columns A:T so we don't
# bloat the sheet dimension out to XFD.
$ws.Range("A31:T31").Copy($ws.Range("A30:T30"))

# Overwrite the new row's unique values (date, volume, min/max/avg price,
# price per kg). The remaining columns (Mercado, Región, Codreg, Tipo,
# Producto, Categoría, Variedad, Calidad, Unidad, Origen, Kg/unidad) keep
# the values copied above.
$ws.Range("D30").Value2 = 45099
$ws.Range("M30").Value2 = 40
$ws.Range("N30").Value2 = 9000
$ws.Range("O30").Value2 = 9000
$ws.Range("P30").Value2 = 9000
$ws.Range("S30").Value2 = 2250
